# Add two new weekly-schedule rows (period 98 and period 99) to the bottom
# of the single worksheet, and move the active selection down to the new
# last row, mirroring the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 48: 2026/2/20 -> 2026/4/17, 第98期 秘寶 開放區域 四象亨通 ...
$ws.Cells.Item(48, 1).Value = "2026/2/20"
$ws.Cells.Item(48, 2).Value = "2026/4/17"
$ws.Cells.Item(48, 3).Value = "第98期 秘寶 開放區域 四象亨通 祕寶效果: 被動進階傷害提高14400000 (34208622)"

# Row 49: 2026/2/27 -> 2026/4/24, 第99期 第八代坐騎
$ws.Cells.Item(49, 1).Value = "2026/2/27"
$ws.Cells.Item(49, 2).Value = "2026/4/24"
$ws.Cells.Item(49, 3).Value = "第99期 第八代坐騎"

# Match the saved selection / scroll position from the diff (D49 selected).
[void]$ws.Range("D49").Select()
